$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices + 1h volume deltas).
# Numeric-looking text values get an explicit Text number format first so
# Excel keeps them as strings (matching the source data), matching the
# original inline-string cell contents instead of being parsed as numbers.

$ws.Range("D2").Value = "28.901.57"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.880.04"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.82"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4605"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3876"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9858"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.819.23"
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.995"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.646"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06957"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.13"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009969"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.99"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "28.897.09"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.235"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.092"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.21"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.008"
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.63"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.923"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09366"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9020"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.262"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.315"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.257"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.637"
$ws.Range("E39").Value = "  -4.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5653"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.672"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.269"
$ws.Range("E43").Value = "  +5.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.93"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07047"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.842"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.064"
$ws.Range("E50").Value = "  -4.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.66"
$ws.Range("E51").Value = "  -0.03%  "

Write-Host "Applied cryptos list update."
